# Append three new data rows (8 / 12 / 16 bit-width runs) to the
# "sim_coeff_selet_gradient_final_" sheet as rows 10-12, and move the
# selection to AB20 (matches the saved workbook's selection state after
# this data was appended).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 10

$data = New-Object 'object[,]' 3,28
$data[0,0] = 8
$data[0,1] = 4
$data[0,2] = 0.125
$data[0,3] = 0.375
$data[0,4] = 1
$data[0,5] = 1
$data[0,6] = 0.5
$data[0,7] = 0.125
$data[0,8] = -0.0625
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 1
$data[0,14] = 2
$data[0,15] = 3
$data[0,16] = 4
$data[0,17] = 4
$data[0,18] = 4
$data[0,19] = 4
$data[0,20] = 100000000000000
$data[0,21] = 0.054495000000000002
$data[0,22] = 0.0077980000000000002
$data[0,23] = 0.0085389999999999997
$data[0,24] = 0.0085679999999999992
$data[0,25] = 0.0085679999999999992
$data[0,26] = 0.0085679999999999992
$data[0,27] = 0.0085679999999999992
$data[1,0] = 12
$data[1,1] = 4
$data[1,2] = 0.125
$data[1,3] = 0.375
$data[1,4] = 1
$data[1,5] = 1
$data[1,6] = 0.5
$data[1,7] = 0.125
$data[1,8] = -0.0625
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 1
$data[1,14] = 2
$data[1,15] = 3
$data[1,16] = 4
$data[1,17] = 4
$data[1,18] = 4
$data[1,19] = 4
$data[1,20] = 100000000000000
$data[1,21] = 0.045599000000000001
$data[1,22] = 0.0045269999999999998
$data[1,23] = 0.00069700000000000003
$data[1,24] = 0.000272
$data[1,25] = 0.000272
$data[1,26] = 0.000272
$data[1,27] = 0.000272
$data[2,0] = 16
$data[2,1] = 4
$data[2,2] = 0.125
$data[2,3] = 0.375
$data[2,4] = 1
$data[2,5] = 1
$data[2,6] = 0.5
$data[2,7] = 0.125
$data[2,8] = -0.0625
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 1
$data[2,14] = 2
$data[2,15] = 3
$data[2,16] = 4
$data[2,17] = 4
$data[2,18] = 4
$data[2,19] = 4
$data[2,20] = 100000000000000
$data[2,21] = 0.045058000000000001
$data[2,22] = 0.0049750000000000003
$data[2,23] = 0.00089899999999999995
$data[2,24] = 0.00027399999999999999
$data[2,25] = 0.00027399999999999999
$data[2,26] = 0.00027399999999999999
$data[2,27] = 0.00027399999999999999

$endRow = $startRow + 2
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 28))
$targetRange.Value2 = $data

$ws.Range("AB20").Select() | Out-Null

Write-Output "done"